$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update example company name in B2
$ws.Range("B2").Value = "(주)서울냉동"

# Fill in the "특이사항" (notes) example value in N2
$ws.Range("N2").Value = "주차공간 협소"

# Widen column N (14) from 6 to 9
# (ColumnWidth uses Excel's character-width units, which the file format
# then re-expresses with its own padding formula; 8.17 here is what lands
# on a stored column width of exactly 9, matching the target file.)
$ws.Columns.Item(14).ColumnWidth = 8.17

# Insert two more blank template rows (rows 4 and 5) below the existing
# blank row 3, pushing nothing else down (they become the new last rows).
$ws.Range("A4:A5").EntireRow.Insert()

# Fill rows 4-5 mirroring row 3's pattern:
# blank A-E, default time window F-I, forklift flag J, duration K, blank L-N
foreach ($r in 4..5) {
    foreach ($col in @(1,2,3,4,5,12,13,14)) {
        $c = $ws.Cells.Item($r, $col)
        $c.Value = "'"
        $c.Style = "Normal"
    }
    $ws.Cells.Item($r, 6).Value = "09:00"
    $ws.Cells.Item($r, 7).Value = "17:00"
    $ws.Cells.Item($r, 8).Value = "09:00"
    $ws.Cells.Item($r, 9).Value = "17:00"
    $ws.Cells.Item($r, 10).Value = "Y"
    $ws.Cells.Item($r, 11).Value = 30
}
